$d = $word.ActiveDocument

# --- Change 1: fix typo "coquillle" -> "coquille" ---
$d.Content.Find.Execute("coquillle de la queue", $true, $false, $false, $false, $false,
                         $true, 1, $false, "coquille de la queue", 2) | Out-Null

# --- Change 2: "ventre entiere<lb/>soict descouverte, pour" becomes
#     "ventre entiere soict<lb/>" + new paragraph + "descouverte, pour<lb/>" ---

# Locate the run containing "soict descouverte, pour"
$find1 = $d.Content
$find1.Find.Execute("soict descouverte, pour", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
$matchStart = $find1.Start

# Remove the leading "soict " (6 characters, including trailing space)
$soictSpace = $d.Range($matchStart, $matchStart + 6)
$soictSpace.Text = ""

# Append " soict" to the end of "ventre entiere"
$find2 = $d.Content
$find2.Find.Execute("ventre entiere", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($find2.End, $find2.End)
$insertPoint.InsertAfter(" soict")

# Split the paragraph right before "descouverte, pour"
$find3 = $d.Content
$find3.Find.Execute("descouverte, pour", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
$breakPoint = $d.Range($find3.Start, $find3.Start)
$breakPoint.InsertParagraphBefore()
